# Weekly fruit/vegetable price update: rows 137-152 shift to make room for a new
# week of "Sandia" (watermelon) price records, and three new rows (153-155) are
# appended so the sheet now spans A1:R155.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 137
$ws.Range("A137").Value = 11
$ws.Range("B137").Value = 'Vega Monumental Concepción'
$ws.Range("C137").Value = 'Bíobío'
$ws.Range("D137").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D137").Value = 44946
$ws.Range("E137").Value = 8
$ws.Range("F137").Value = 100112028
$ws.Range("G137").Value = 'Sandia'
$ws.Range("H137").Value = 'Sin especificar'
$ws.Range("I137").Value = 'Primera'
$ws.Range("J137").Value = 1000
$ws.Range("K137").Value = 3000
$ws.Range("L137").Value = 3000
$ws.Range("M137").Value = 3000
$ws.Range("N137").Value = '$/unidad'
$ws.Range("O137").Value = 'Región de O''Higgins'
$ws.Range("P137").Value = 3000
$ws.Range("Q137").Value = 1
$ws.Range("R137").Value = 'Hortaliza'

# Row 138
$ws.Range("A138").Value = 11
$ws.Range("B138").Value = 'Vega Monumental Concepción'
$ws.Range("C138").Value = 'Bíobío'
$ws.Range("D138").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D138").Value = 44946
$ws.Range("E138").Value = 8
$ws.Range("F138").Value = 100112028
$ws.Range("G138").Value = 'Sandia'
$ws.Range("H138").Value = 'Sin especificar'
$ws.Range("I138").Value = 'Segunda'
$ws.Range("J138").Value = 1000
$ws.Range("K138").Value = 2300
$ws.Range("L138").Value = 2300
$ws.Range("M138").Value = 2300
$ws.Range("N138").Value = '$/unidad'
$ws.Range("O138").Value = 'Región de O''Higgins'
$ws.Range("P138").Value = 2300
$ws.Range("Q138").Value = 1
$ws.Range("R138").Value = 'Hortaliza'

# Row 139
$ws.Range("A139").Value = 11
$ws.Range("B139").Value = 'Vega Monumental Concepción'
$ws.Range("C139").Value = 'Bíobío'
$ws.Range("D139").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D139").Value = 44946
$ws.Range("E139").Value = 8
$ws.Range("F139").Value = 100112028
$ws.Range("G139").Value = 'Sandia'
$ws.Range("H139").Value = 'Sin especificar'
$ws.Range("I139").Value = 'Tercera'
$ws.Range("J139").Value = 1000
$ws.Range("K139").Value = 2000
$ws.Range("L139").Value = 2000
$ws.Range("M139").Value = 2000
$ws.Range("N139").Value = '$/unidad'
$ws.Range("O139").Value = 'Región de O''Higgins'
$ws.Range("P139").Value = 2000
$ws.Range("Q139").Value = 1
$ws.Range("R139").Value = 'Hortaliza'

# Row 140
$ws.Range("A140").Value = 11
$ws.Range("B140").Value = 'Vega Monumental Concepción'
$ws.Range("C140").Value = 'Bíobío'
$ws.Range("D140").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D140").Value = 44581
$ws.Range("E140").Value = 8
$ws.Range("F140").Value = 100112028
$ws.Range("G140").Value = 'Sandia'
$ws.Range("H140").Value = 'Sin especificar'
$ws.Range("I140").Value = 'Extra'
$ws.Range("J140").Value = 400
$ws.Range("K140").Value = 2500
$ws.Range("L140").Value = 2500
$ws.Range("M140").Value = 2500
$ws.Range("N140").Value = '$/unidad'
$ws.Range("O140").Value = 'Región del Maule'
$ws.Range("P140").Value = 2500
$ws.Range("Q140").Value = 1
$ws.Range("R140").Value = 'Hortaliza'

# Row 141
$ws.Range("A141").Value = 11
$ws.Range("B141").Value = 'Vega Monumental Concepción'
$ws.Range("C141").Value = 'Bíobío'
$ws.Range("D141").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D141").Value = 44581
$ws.Range("E141").Value = 8
$ws.Range("F141").Value = 100112028
$ws.Range("G141").Value = 'Sandia'
$ws.Range("H141").Value = 'Sin especificar'
$ws.Range("I141").Value = 'Primera'
$ws.Range("J141").Value = 400
$ws.Range("K141").Value = 2000
$ws.Range("L141").Value = 2000
$ws.Range("M141").Value = 2000
$ws.Range("N141").Value = '$/unidad'
$ws.Range("O141").Value = 'Región del Maule'
$ws.Range("P141").Value = 2000
$ws.Range("Q141").Value = 1
$ws.Range("R141").Value = 'Hortaliza'

# Row 142
$ws.Range("A142").Value = 11
$ws.Range("B142").Value = 'Vega Monumental Concepción'
$ws.Range("C142").Value = 'Bíobío'
$ws.Range("D142").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D142").Value = 44581
$ws.Range("E142").Value = 8
$ws.Range("F142").Value = 100112028
$ws.Range("G142").Value = 'Sandia'
$ws.Range("H142").Value = 'Sin especificar'
$ws.Range("I142").Value = 'Segunda'
$ws.Range("J142").Value = 400
$ws.Range("K142").Value = 1500
$ws.Range("L142").Value = 1500
$ws.Range("M142").Value = 1500
$ws.Range("N142").Value = '$/unidad'
$ws.Range("O142").Value = 'Región del Maule'
$ws.Range("P142").Value = 1500
$ws.Range("Q142").Value = 1
$ws.Range("R142").Value = 'Hortaliza'

# Row 143
$ws.Range("A143").Value = 11
$ws.Range("B143").Value = 'Vega Monumental Concepción'
$ws.Range("C143").Value = 'Bíobío'
$ws.Range("D143").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D143").Value = 44908
$ws.Range("E143").Value = 8
$ws.Range("F143").Value = 100112028
$ws.Range("G143").Value = 'Sandia'
$ws.Range("H143").Value = 'Sin especificar'
$ws.Range("I143").Value = 'Extra'
$ws.Range("J143").Value = 500
$ws.Range("K143").Value = 4000
$ws.Range("L143").Value = 4000
$ws.Range("M143").Value = 4000
$ws.Range("N143").Value = '$/unidad'
$ws.Range("O143").Value = 'Región de O''Higgins'
$ws.Range("P143").Value = 4000
$ws.Range("Q143").Value = 1
$ws.Range("R143").Value = 'Hortaliza'

# Row 144
$ws.Range("A144").Value = 11
$ws.Range("B144").Value = 'Vega Monumental Concepción'
$ws.Range("C144").Value = 'Bíobío'
$ws.Range("D144").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D144").Value = 44908
$ws.Range("E144").Value = 8
$ws.Range("F144").Value = 100112028
$ws.Range("G144").Value = 'Sandia'
$ws.Range("H144").Value = 'Sin especificar'
$ws.Range("I144").Value = 'Primera'
$ws.Range("J144").Value = 500
$ws.Range("K144").Value = 3500
$ws.Range("L144").Value = 3500
$ws.Range("M144").Value = 3500
$ws.Range("N144").Value = '$/unidad'
$ws.Range("O144").Value = 'Región de O''Higgins'
$ws.Range("P144").Value = 3500
$ws.Range("Q144").Value = 1
$ws.Range("R144").Value = 'Hortaliza'

# Row 145
$ws.Range("A145").Value = 11
$ws.Range("B145").Value = 'Vega Monumental Concepción'
$ws.Range("C145").Value = 'Bíobío'
$ws.Range("D145").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D145").Value = 44908
$ws.Range("E145").Value = 8
$ws.Range("F145").Value = 100112028
$ws.Range("G145").Value = 'Sandia'
$ws.Range("H145").Value = 'Sin especificar'
$ws.Range("I145").Value = 'Segunda'
$ws.Range("J145").Value = 500
$ws.Range("K145").Value = 3000
$ws.Range("L145").Value = 3000
$ws.Range("M145").Value = 3000
$ws.Range("N145").Value = '$/unidad'
$ws.Range("O145").Value = 'Región de O''Higgins'
$ws.Range("P145").Value = 3000
$ws.Range("Q145").Value = 1
$ws.Range("R145").Value = 'Hortaliza'

# Row 146
$ws.Range("A146").Value = 11
$ws.Range("B146").Value = 'Vega Monumental Concepción'
$ws.Range("C146").Value = 'Bíobío'
$ws.Range("D146").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D146").Value = 44281
$ws.Range("E146").Value = 8
$ws.Range("F146").Value = 100112028
$ws.Range("G146").Value = 'Sandia'
$ws.Range("H146").Value = 'Sin especificar'
$ws.Range("I146").Value = 'Primera'
$ws.Range("J146").Value = 400
$ws.Range("K146").Value = 2500
$ws.Range("L146").Value = 2500
$ws.Range("M146").Value = 2500
$ws.Range("N146").Value = '$/unidad'
$ws.Range("O146").Value = 'Región de O''Higgins'
$ws.Range("P146").Value = 2500
$ws.Range("Q146").Value = 1
$ws.Range("R146").Value = 'Hortaliza'

# Row 147
$ws.Range("A147").Value = 11
$ws.Range("B147").Value = 'Vega Monumental Concepción'
$ws.Range("C147").Value = 'Bíobío'
$ws.Range("D147").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D147").Value = 44281
$ws.Range("E147").Value = 8
$ws.Range("F147").Value = 100112028
$ws.Range("G147").Value = 'Sandia'
$ws.Range("H147").Value = 'Sin especificar'
$ws.Range("I147").Value = 'Segunda'
$ws.Range("J147").Value = 400
$ws.Range("K147").Value = 2000
$ws.Range("L147").Value = 2000
$ws.Range("M147").Value = 2000
$ws.Range("N147").Value = '$/unidad'
$ws.Range("O147").Value = 'Región de O''Higgins'
$ws.Range("P147").Value = 2000
$ws.Range("Q147").Value = 1
$ws.Range("R147").Value = 'Hortaliza'

# Row 148
$ws.Range("A148").Value = 11
$ws.Range("B148").Value = 'Vega Monumental Concepción'
$ws.Range("C148").Value = 'Bíobío'
$ws.Range("D148").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D148").Value = 44251
$ws.Range("E148").Value = 8
$ws.Range("F148").Value = 100112028
$ws.Range("G148").Value = 'Sandia'
$ws.Range("H148").Value = 'Sin especificar'
$ws.Range("I148").Value = 'Extra'
$ws.Range("J148").Value = 300
$ws.Range("K148").Value = 3000
$ws.Range("L148").Value = 3000
$ws.Range("M148").Value = 3000
$ws.Range("N148").Value = '$/unidad'
$ws.Range("O148").Value = 'Región de O''Higgins'
$ws.Range("P148").Value = 3000
$ws.Range("Q148").Value = 1
$ws.Range("R148").Value = 'Hortaliza'

# Row 149
$ws.Range("A149").Value = 11
$ws.Range("B149").Value = 'Vega Monumental Concepción'
$ws.Range("C149").Value = 'Bíobío'
$ws.Range("D149").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D149").Value = 44251
$ws.Range("E149").Value = 8
$ws.Range("F149").Value = 100112028
$ws.Range("G149").Value = 'Sandia'
$ws.Range("H149").Value = 'Sin especificar'
$ws.Range("I149").Value = 'Primera'
$ws.Range("J149").Value = 500
$ws.Range("K149").Value = 2500
$ws.Range("L149").Value = 2500
$ws.Range("M149").Value = 2500
$ws.Range("N149").Value = '$/unidad'
$ws.Range("O149").Value = 'Región de O''Higgins'
$ws.Range("P149").Value = 2500
$ws.Range("Q149").Value = 1
$ws.Range("R149").Value = 'Hortaliza'

# Row 150
$ws.Range("A150").Value = 11
$ws.Range("B150").Value = 'Vega Monumental Concepción'
$ws.Range("C150").Value = 'Bíobío'
$ws.Range("D150").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D150").Value = 44251
$ws.Range("E150").Value = 8
$ws.Range("F150").Value = 100112028
$ws.Range("G150").Value = 'Sandia'
$ws.Range("H150").Value = 'Sin especificar'
$ws.Range("I150").Value = 'Segunda'
$ws.Range("J150").Value = 500
$ws.Range("K150").Value = 2000
$ws.Range("L150").Value = 2000
$ws.Range("M150").Value = 2000
$ws.Range("N150").Value = '$/unidad'
$ws.Range("O150").Value = 'Región de O''Higgins'
$ws.Range("P150").Value = 2000
$ws.Range("Q150").Value = 1
$ws.Range("R150").Value = 'Hortaliza'

# Row 151
$ws.Range("A151").Value = 11
$ws.Range("B151").Value = 'Vega Monumental Concepción'
$ws.Range("C151").Value = 'Bíobío'
$ws.Range("D151").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D151").Value = 44225
$ws.Range("E151").Value = 8
$ws.Range("F151").Value = 100112028
$ws.Range("G151").Value = 'Sandia'
$ws.Range("H151").Value = 'Sin especificar'
$ws.Range("I151").Value = 'Extra'
$ws.Range("J151").Value = 500
$ws.Range("K151").Value = 3000
$ws.Range("L151").Value = 3000
$ws.Range("M151").Value = 3000
$ws.Range("N151").Value = '$/unidad'
$ws.Range("O151").Value = 'Región de O''Higgins'
$ws.Range("P151").Value = 3000
$ws.Range("Q151").Value = 1
$ws.Range("R151").Value = 'Hortaliza'

# Row 152
$ws.Range("A152").Value = 11
$ws.Range("B152").Value = 'Vega Monumental Concepción'
$ws.Range("C152").Value = 'Bíobío'
$ws.Range("D152").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D152").Value = 44225
$ws.Range("E152").Value = 8
$ws.Range("F152").Value = 100112028
$ws.Range("G152").Value = 'Sandia'
$ws.Range("H152").Value = 'Sin especificar'
$ws.Range("I152").Value = 'Primera'
$ws.Range("J152").Value = 500
$ws.Range("K152").Value = 2500
$ws.Range("L152").Value = 2500
$ws.Range("M152").Value = 2500
$ws.Range("N152").Value = '$/unidad'
$ws.Range("O152").Value = 'Región de O''Higgins'
$ws.Range("P152").Value = 2500
$ws.Range("Q152").Value = 1
$ws.Range("R152").Value = 'Hortaliza'

# Row 153
$ws.Range("A153").Value = 11
$ws.Range("B153").Value = 'Vega Monumental Concepción'
$ws.Range("C153").Value = 'Bíobío'
$ws.Range("D153").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D153").Value = 44225
$ws.Range("E153").Value = 8
$ws.Range("F153").Value = 100112028
$ws.Range("G153").Value = 'Sandia'
$ws.Range("H153").Value = 'Sin especificar'
$ws.Range("I153").Value = 'Segunda'
$ws.Range("J153").Value = 500
$ws.Range("K153").Value = 2000
$ws.Range("L153").Value = 2000
$ws.Range("M153").Value = 2000
$ws.Range("N153").Value = '$/unidad'
$ws.Range("O153").Value = 'Región de O''Higgins'
$ws.Range("P153").Value = 2000
$ws.Range("Q153").Value = 1
$ws.Range("R153").Value = 'Hortaliza'

# Row 154
$ws.Range("A154").Value = 11
$ws.Range("B154").Value = 'Vega Monumental Concepción'
$ws.Range("C154").Value = 'Bíobío'
$ws.Range("D154").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D154").Value = 44272
$ws.Range("E154").Value = 8
$ws.Range("F154").Value = 100112028
$ws.Range("G154").Value = 'Sandia'
$ws.Range("H154").Value = 'Sin especificar'
$ws.Range("I154").Value = 'Primera'
$ws.Range("J154").Value = 300
$ws.Range("K154").Value = 2500
$ws.Range("L154").Value = 2500
$ws.Range("M154").Value = 2500
$ws.Range("N154").Value = '$/unidad'
$ws.Range("O154").Value = 'Región de O''Higgins'
$ws.Range("P154").Value = 2500
$ws.Range("Q154").Value = 1
$ws.Range("R154").Value = 'Hortaliza'

# Row 155
$ws.Range("A155").Value = 11
$ws.Range("B155").Value = 'Vega Monumental Concepción'
$ws.Range("C155").Value = 'Bíobío'
$ws.Range("D155").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D155").Value = 44272
$ws.Range("E155").Value = 8
$ws.Range("F155").Value = 100112028
$ws.Range("G155").Value = 'Sandia'
$ws.Range("H155").Value = 'Sin especificar'
$ws.Range("I155").Value = 'Segunda'
$ws.Range("J155").Value = 300
$ws.Range("K155").Value = 2000
$ws.Range("L155").Value = 2000
$ws.Range("M155").Value = 2000
$ws.Range("N155").Value = '$/unidad'
$ws.Range("O155").Value = 'Región de O''Higgins'
$ws.Range("P155").Value = 2000
$ws.Range("Q155").Value = 1
$ws.Range("R155").Value = 'Hortaliza'

